$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
$wsMeta.Range("B8").Value = "2025-05-05T11:54:16+00:00"
$wsMeta.Range("B17").Value = "string"
$wsMeta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/string"

# --- Elements sheet updates (integer -> string primitive type) ---

# Row 2: root element "integer" -> "string"
$wsElem.Range("A2").Value = "string"
$wsElem.Range("B2").Value = "string"
$wsElem.Range("L2").Value = "Primitive Type string"
$wsElem.Range("M2").Value = "A sequence of Unicode characters"
$wsElem.Range("N2").Value = "Note that FHIR strings SHALL NOT exceed 1MB in size"
$wsElem.Range("AF2").Value = "string"

# Row 3: "integer.id" -> "string.id"
$wsElem.Range("A3").Value = "string.id"
$wsElem.Range("B3").Value = "string.id"

# Row 4: "integer.extension" -> "string.extension"
$wsElem.Range("A4").Value = "string.extension"
$wsElem.Range("B4").Value = "string.extension"

# Row 5: "integer.value" -> "string.value"
$wsElem.Range("A5").Value = "string.value"
$wsElem.Range("B5").Value = "string.value"
$wsElem.Range("K5").Value = "string`n"
$wsElem.Range("L5").Value = "Primitive value for string"
$wsElem.Range("M5").Value = "The actual value"
# Minimum/Maximum Value no longer apply to string -> cleared (kept as text cells)
$wsElem.Range("U5").Value = "'"
$wsElem.Range("V5").Value = "'"
# Maximum Length is stored as text in this sheet (like the other numeric-looking
# Min/Max columns), so force a leading-apostrophe to keep it a text value.
$wsElem.Range("W5").Value = "'1048576"
$wsElem.Range("AF5").Value = "string.value"

# Column A/B best-fit width shrinks now that "integer.extension" (17 chars)
# became "string.extension" (16 chars) -- the longest string in those columns.
$wsElem.Columns.Item(1).AutoFit() | Out-Null
$wsElem.Columns.Item(2).AutoFit() | Out-Null
